# Apply the "Updated cryptos list" GitHub Actions refresh to the crypto
# tracker worksheet: row 39/40 coin entries are swapped (RenderToken now
# ranks above FirstDigitalUSD) and the Price / Volume(1h) columns are
# refreshed with newly scraped values across the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 39 / 40 coins swapped places, each with refreshed price & volume ---
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D39").Formula = "'6.21"
$ws.Range("E39").Formula = "'  +11.52%  "

$ws.Range("B40").Value = "FirstDigitalUSD"
$ws.Range("C40").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D40").Formula = "'0.999"
$ws.Range("E40").Formula = "'  -0.02%  "

# --- Refresh Price (D) and Volume(1h) (E) columns for all other rows ---
$ws.Range("D2").Formula = "'82.364.45"
$ws.Range("E2").Formula = "'  +3.63%  "
$ws.Range("D3").Formula = "'3.181.18"
$ws.Range("E3").Formula = "'  -0.10%  "
$ws.Range("D4").Formula = "'0.999"
$ws.Range("E4").Formula = "'  -0.21%  "
$ws.Range("D5").Formula = "'219.68"
$ws.Range("E5").Formula = "'  +6.93%  "
$ws.Range("D6").Formula = "'619.50"
$ws.Range("E6").Formula = "'  -2.71%  "
$ws.Range("D7").Formula = "'0.290"
$ws.Range("E7").Formula = "'  +20.00%  "
$ws.Range("D8").Formula = "'0.999"
$ws.Range("E8").Formula = "'  -0.07%  "
$ws.Range("D9").Formula = "'0.582"
$ws.Range("E9").Formula = "'  -3.43%  "
$ws.Range("D10").Formula = "'3.173.17"
$ws.Range("E10").Formula = "'  -0.44%  "
$ws.Range("D11").Formula = "'0.592"
$ws.Range("E11").Formula = "'  -2.07%  "
$ws.Range("D12").Formula = "'0.0000256"
$ws.Range("E12").Formula = "'  -0.30%  "
$ws.Range("E13").Formula = "'  -0.57%  "
$ws.Range("D14").Formula = "'5.31"
$ws.Range("E14").Formula = "'  -1.33%  "
$ws.Range("D15").Formula = "'3.747.37"
$ws.Range("E15").Formula = "'  -0.93%  "
$ws.Range("D16").Formula = "'32.21"
$ws.Range("E16").Formula = "'  -0.10%  "
$ws.Range("D17").Formula = "'81.995.20"
$ws.Range("E17").Formula = "'  +3.03%  "
$ws.Range("D18").Formula = "'3.163.60"
$ws.Range("E18").Formula = "'  -1.40%  "
$ws.Range("D19").Formula = "'3.24"
$ws.Range("E19").Formula = "'  +10.31%  "
$ws.Range("D20").Formula = "'13.97"
$ws.Range("E20").Formula = "'  -4.10%  "
$ws.Range("D21").Formula = "'438.27"
$ws.Range("E21").Formula = "'  -0.12%  "
$ws.Range("D22").Formula = "'8.90"
$ws.Range("E22").Formula = "'  -5.53%  "
$ws.Range("D23").Formula = "'5.13"
$ws.Range("E23").Formula = "'  -2.11%  "
$ws.Range("D24").Formula = "'7.29"
$ws.Range("E24").Formula = "'  +4.67%  "
$ws.Range("D25").Formula = "'5.25"
$ws.Range("E25").Formula = "'  +9.47%  "
$ws.Range("D26").Formula = "'11.93"
$ws.Range("E26").Formula = "'  +10.11%  "
$ws.Range("E27").Formula = "'  -1.60%  "
$ws.Range("D28").Formula = "'76.65"
$ws.Range("E28").Formula = "'  -1.25%  "
$ws.Range("D29").Formula = "'0.999"
$ws.Range("E29").Formula = "'  -0.18%  "
$ws.Range("D30").Formula = "'0.0000121"
$ws.Range("E30").Formula = "'  -2.11%  "
$ws.Range("D31").Formula = "'0.998"
$ws.Range("E31").Formula = "'  -0.27%  "
$ws.Range("D32").Formula = "'9.02"
$ws.Range("E32").Formula = "'  -1.06%  "
$ws.Range("D33").Formula = "'569.41"
$ws.Range("E33").Formula = "'  +4.63%  "
$ws.Range("D34").Formula = "'1.48"
$ws.Range("E34").Formula = "'  -4.24%  "
$ws.Range("D35").Formula = "'0.147"
$ws.Range("E35").Formula = "'  +21.18%  "
$ws.Range("D36").Formula = "'0.151"
$ws.Range("E36").Formula = "'  -0.59%  "
$ws.Range("D37").Formula = "'1.98"
$ws.Range("E37").Formula = "'  -2.70%  "
$ws.Range("D38").Formula = "'22.62"
$ws.Range("E38").Formula = "'  -1.91%  "
$ws.Range("D41").Formula = "'0.404"
$ws.Range("E41").Formula = "'  -2.23%  "
$ws.Range("D42").Formula = "'20.84"
$ws.Range("E42").Formula = "'  +4.05%  "
$ws.Range("D43").Formula = "'2.01"
$ws.Range("E43").Formula = "'  +11.60%  "
$ws.Range("D44").Formula = "'3.01"
$ws.Range("E44").Formula = "'  +12.72%  "
$ws.Range("D45").Formula = "'159.23"
$ws.Range("E45").Formula = "'  -3.23%  "
$ws.Range("E46").Formula = "'  +0.00%  "
$ws.Range("D47").Formula = "'186.21"
$ws.Range("E47").Formula = "'  -3.53%  "
$ws.Range("D48").Formula = "'44.67"
$ws.Range("E48").Formula = "'  +1.59%  "
$ws.Range("D49").Formula = "'1.32"
$ws.Range("E49").Formula = "'  -0.86%  "
$ws.Range("D50").Formula = "'0.768"
$ws.Range("E50").Formula = "'  -3.94%  "
$ws.Range("D51").Formula = "'25.83"
$ws.Range("E51").Formula = "'  -0.55%  "
